# Testimony.docx edit — merge the two paragraphs into one and rewrite
# the 'meeting Christ' narrative per the commit message / diff.
$d = $word.ActiveDocument

# --- Step 1: merge paragraph 1 and paragraph 2 into a single paragraph ---
# (the edit drops the paragraph break after '...I suffered.')
$p1 = $d.Paragraphs(1)
$markRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$markRange.Delete()

# --- Step 2: text replacement #1 ---
$old = @'
My struggle has always been with relationships. In high school, I invested my happiness into my friendships. This meant that my happiest moments were when my friends and I were together, laughing over dumb jokes or playing video games together. This also means that my lowest moments were when they left, and I was left to scrape the internet for any means of social connection. Thus, I loved my friends, but I loved them because they filled that hole; they made me feel like more than nothing. The person who suffered the most from this, besides me, was my best friend Ivy. We were constant friends, and she was the person I would always message if I felt alone or empty. For that, I loved her, but this was not a Godly love. During the summer before college, she became so busy with her job that we couldn’t talk regularly. This is where my love fell short. I became so embittered with her because she wouldn’t fill my loneliness any more, and I wanted her to suffer as deeply as I suffered. 
'@
$old = $old.TrimEnd("`r","`n")
$new = @'
In high school, I had a lot of friends of convenience: math friends, track friends, band friends, etc. The only friend who wasn’t a friend of convenience was my best friend Ivy. We were constants in each other’s lives for the whole of high school. We were each other’s comforts through bad tests, hard classes, and messy relationships. I loved Ivy because often, she was what I needed to get through the day. Whenever I felt alone, I would message her with whatever was on my mind. It was a relationship where I was looking out for myself. During the summer before college, she became so busy with her job that we couldn’t talk regularly. When she couldn’t make me feel less lonely, my love for her fell short. I became so embittered with her that I didn’t talk to her for a year because I wanted her to suffer as deeply as I suffered.
'@
$new = $new.TrimEnd("`r","`n")
$found0 = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Replacement 1 applied: $found0"

# --- Step 3: text replacement #2 ---
$old = @'
was so desperately looking for. I had never met people who were so willing to sacrifice for you, whether it was their time or their money. 
'@
$old = $old.TrimEnd("`r","`n")
$new = @'
was so desperately looking for. 
'@
$new = $new.TrimEnd("`r","`n")
$found1 = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Replacement 2 applied: $found1"

# --- Step 4: text replacement #3 ---
$old = @'
, but also because I wanted their love to flow into me. I went to services, and listened to sermons about a God who loved, and in His love, I saw how imperfect my love is. I saw how I withheld it from people I thought didn’t deserve it, and how I loved for myself. The clearest place where I was withholding love was with Ivy, the best friend on whom I turned my back, and I knew that God would want me to forgive her, but I couldn’t. There wasn’t enough of God’s love in my heart to overcome my love for myself. 
'@
$old = $old.TrimEnd("`r","`n")
$new = @'
, but also because I saw their love and I wanted more. Through family group and sermons, I learned about a God who loved, and in His love, I saw how imperfect my love is. Where I only loved “deserving” people, God loved indiscriminately. Where I loved because of what people did, God loved because of who we are. Where I turned away from the people that sinned against me, God sent His son down so that we could be reconciled with Him. There wasn’t enough of God’s love in my heart to overcome my love for myself. 
'@
$new = $new.TrimEnd("`r","`n")
$found2 = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Output "Replacement 3 applied: $found2"

# --- Step 5: reposition the '_GoBack' last-edit bookmark ---
# In the final text, the last block of newly authored material ends right
# after '...reconciled with Him. ' — that's where Word would drop _GoBack.
$anchor = $d.Content
$anchorFound = $anchor.Find.Execute("reconciled with Him. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($anchorFound) {
    $goBackPoint = $d.Range($anchor.End, $anchor.End)
    $d.Bookmarks.Add("_GoBack", $goBackPoint)
}
Write-Output "_GoBack repositioned: $anchorFound"

Write-Output $d.Content.Text
